# Fixed a bug in symbolTrigger
#
# The symbol/reel weight table on Sheet1 (A2:F25) was generated in the
# wrong row order. This restores the correct row ordering for the
# weighted-symbol table (column A = symbol id, columns B:F = reel1..reel5
# weights) without touching the totals row (row 26) or any rows whose
# position was already correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $a, $b, $c, $d, $e, $f) {
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
}

Set-Row  2  601   9  60  67  60  42
Set-Row  4  701   3  90  45  97  15
Set-Row  5  401   9  48  67  75  45
Set-Row  6  801   3  67  65  52  45
Set-Row  7  101   9  30  15  60  15
Set-Row  8  501   9  52  30  75  45
Set-Row  9  1202  2  10  10  10  10
Set-Row 10  1001 18  30  75  60  72
Set-Row 11  902   1   0   0   0   0
Set-Row 12  301   6  45  30  60  45
Set-Row 13  201   9  30  15  45  30
Set-Row 14  1201  2  10  10  10  10
Set-Row 15  901  16  15  45  60  60
Set-Row 16  1    0   2   2   2   2
Set-Row 17  1101  0  15  30  30   0
Set-Row 19  2    0   2   2   2   2
Set-Row 20  802   0   4   5   4   0
Set-Row 22  402   0   0   4   0   0
Set-Row 23  602   0   0   4   0   9
